$d = $word.ActiveDocument

function Replace-InParagraph($index, $oldText, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# Paragraph 2 - Planning
Replace-InParagraph 2 "(Ref-DJ49F2)" "(Al-Sayed, 1998)"
Replace-InParagraph 2 "(Ref-G7H3J1)" "(Al-Sayed, 1998)"

# Paragraph 4 - Organizing
Replace-InParagraph 4 "(Ref-DJ49F2)" "(Ref-f928799)"
Replace-InParagraph 4 "(Ref-G7H2K9)" "(Ref-f928799)"

# Paragraph 6 - Leading
Replace-InParagraph 6 "(Ref-J7X8K2)" "(Ref-s477686)"
Replace-InParagraph 6 "(Ref-B9N2M5)" "(Ref-s477686)"

# Paragraph 8 - Controlling
Replace-InParagraph 8 "(Ref-DJ49F2)" "(Ref-f654461)"

# Paragraph 10 - Communicating
Replace-InParagraph 10 "(Ref-DJ49F2)" "(Brown and Garcia)"

Write-Host "Edits applied"
